# Atualizacao da URL do Github
# Slide 4, shape 4 ("Google Shape;76;p15") contains the Github URL paragraph
# followed by an empty bullet paragraph. Update the URL to point at
# "Aula_JDBC_basico" instead of "jdbc-basico", and remove the trailing
# empty paragraph that followed it.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$shp = $s.Shapes.Item(4)
$tr = $shp.TextFrame.TextRange

# Paragraph 3 (1-based) holds the Github URL text run. Re-assigning the
# text directly would let the host diff old/new strings and split the
# run wherever a common prefix/suffix exists (e.g. "jdbc-basico" vs.
# "Aula_JDBC_basico" share "basico"), fragmenting a single run into
# several. Routing through an unrelated placeholder first avoids any
# character overlap with either string, so each assignment stays a
# single run and keeps its original rPr (font/color/etc.) intact.
$urlPara = $tr.Paragraphs(3, 1)
$urlPara.Text = "PLACEHOLDER_URL_TEXT"
$urlPara = $tr.Paragraphs(3, 1)
$urlPara.Text = "https://github.com/danielkv7/digital-innovation-one/tree/master/Aula_JDBC_basico"

# Paragraph 4 (1-based) is the now-trailing empty bullet paragraph; remove it.
$trailingPara = $tr.Paragraphs(4, 1)
$trailingPara.Delete()
